$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.190.48'
$ws.Range('E2').Value = '  -3.51%  '

$ws.Range('D3').Value = '3.301.22'
$ws.Range('E3').Value = '  +0.65%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '''569.70'
$ws.Range('E5').Value = '  -2.77%  '

$ws.Range('D6').Value = '''180.44'
$ws.Range('E6').Value = '  -3.94%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').Value = '''0.598'
$ws.Range('E8').Value = '  -0.70%  '

$ws.Range('D9').Value = '3.299.68'
$ws.Range('E9').Value = '  +0.62%  '

$ws.Range('E10').Value = '  -1.74%  '

$ws.Range('D11').Value = '''6.63'
$ws.Range('E11').Value = '  -0.27%  '

$ws.Range('D12').Value = '''0.402'
$ws.Range('E12').Value = '  -2.51%  '

$ws.Range('D13').Value = '3.880.04'
$ws.Range('E13').Value = '  +0.65%  '

$ws.Range('E14').Value = '  -0.84%  '

$ws.Range('D15').Value = '''26.90'
$ws.Range('E15').Value = '  -2.63%  '

$ws.Range('D16').Value = '66.311.22'
$ws.Range('E16').Value = '  -3.36%  '

$ws.Range('D17').Value = '''0.0000166'
$ws.Range('E17').Value = '  -2.08%  '

$ws.Range('D18').Value = '3.331.69'
$ws.Range('E18').Value = '  +1.55%  '

$ws.Range('D19').Value = '''438.23'
$ws.Range('E19').Value = '  +5.37%  '

$ws.Range('D20').Value = '''13.55'
$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('D21').Value = '''5.66'
$ws.Range('E21').Value = '  -1.56%  '

$ws.Range('D22').Value = '''7.64'
$ws.Range('E22').Value = '  +0.77%  '

$ws.Range('D23').Value = '''73.94'
$ws.Range('E23').Value = '  +3.32%  '

$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  -0.12%  '

$ws.Range('D25').Value = '3.464.93'
$ws.Range('E25').Value = '  +1.07%  '

$ws.Range('D26').Value = '''0.512'
$ws.Range('E26').Value = '  +0.26%  '

$ws.Range('D27').Value = '''0.0000118'
$ws.Range('E27').Value = '  -0.67%  '

$ws.Range('D28').Value = '''0.191'
$ws.Range('E28').Value = '  +0.27%  '

$ws.Range('D29').Value = '''8.93'
$ws.Range('E29').Value = '  -5.91%  '

$ws.Range('E30').Value = '  +0.02%  '

$ws.Range('D31').Value = '''1.95'
$ws.Range('E31').Value = '  -0.50%  '

$ws.Range('D32').Value = '''22.76'
$ws.Range('E32').Value = '  -0.03%  '

$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = '''0.998'
$ws.Range('E33').Value = '  -0.04%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''5.28'
$ws.Range('E34').Value = '  -4.27%  '

$ws.Range('D35').Value = '''6.76'
$ws.Range('E35').Value = '  -2.05%  '

$ws.Range('D36').Value = '''1.21'
$ws.Range('E36').Value = '  -3.56%  '

$ws.Range('D37').Value = '''1.49'
$ws.Range('E37').Value = '  +1.45%  '

$ws.Range('D38').Value = '''159.97'
$ws.Range('E38').Value = '  -2.60%  '

$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '''27.20'
$ws.Range('E39').Value = '  +2.46%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''1.84'
$ws.Range('E40').Value = '  -3.82%  '

$ws.Range('D41').Value = '2.796.45'
$ws.Range('E41').Value = '  +4.24%  '

$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '''0.784'
$ws.Range('E42').Value = '  -1.37%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '''4.45'
$ws.Range('E43').Value = '  -1.87%  '

$ws.Range('D44').Value = '''6.19'
$ws.Range('E44').Value = '  -3.55%  '

$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '''0.0673'
$ws.Range('E45').Value = '  -1.46%  '

$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '''40.20'
$ws.Range('E46').Value = '  -0.88%  '

$ws.Range('D47').Value = '''24.10'
$ws.Range('E47').Value = '  -2.69%  '

$ws.Range('D48').Value = '''2.31'
$ws.Range('E48').Value = '  -5.56%  '

$ws.Range('D49').Value = '''317.62'
$ws.Range('E49').Value = '  -6.03%  '

$ws.Range('D50').Value = '''0.0270'
$ws.Range('E50').Value = '  -2.24%  '

$ws.Range('D51').Value = '''0.976'
$ws.Range('E51').Value = '  -1.96%  '
